$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.404.48'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '1.822.10'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''314.31'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D7').Value = '''0.5123'
$ws.Range('E7').Value = '  -3.43%  '
$ws.Range('D8').Value = '''0.3930'
$ws.Range('E8').Value = '  -2.86%  '
$ws.Range('D9').Value = '''0.07658'
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '''41.65'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '''1.107'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').Value = '''20.98'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '''6.268'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').Value = '''1.002'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '''7.489'
$ws.Range('E15').Value = '  -1.45%  '
$ws.Range('D16').Value = '1.823.87'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = '''93.20'
$ws.Range('E17').Value = '  +4.00%  '
$ws.Range('D18').Value = '''0.00001097'
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('D19').Value = '''0.06648'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').Value = '''17.72'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '''6.102'
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('D23').Value = '28.423.23'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '''11.16'
$ws.Range('E24').Value = '  -1.45%  '
$ws.Range('D25').Value = '''2.255'
$ws.Range('E25').Value = '  +6.48%  '
$ws.Range('D26').Value = '''20.80'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('D27').Value = '''156.04'
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('D28').Value = '2.033.86'
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D29').Value = '''2.388'
$ws.Range('E29').Value = '  -3.50%  '
$ws.Range('D30').Value = '''123.91'
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').Value = '''1.107'
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('D32').Value = '''0.1095'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').Value = '''5.650'
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('D34').Value = '''3.657'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = '''0.07075'
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('D37').Value = '''0.02325'
$ws.Range('E37').Value = '  -1.04%  '
$ws.Range('D38').Value = '''5.167'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('D39').Value = '''8.750'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').Value = '''0.6249'
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('D41').Value = '''11.19'
$ws.Range('E41').Value = '  -1.49%  '
$ws.Range('D42').Value = '''1.172'
$ws.Range('E42').Value = '  -1.80%  '
$ws.Range('D44').Value = '''1.392'
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('D45').Value = '''13.37'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('D47').Value = '''0.5877'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').Value = '''124.50'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').Value = '''1.978'
$ws.Range('D50').Value = '''1.194'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = '''0.06896'
$ws.Range('E51').Value = '  -0.15%  '
